$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031756065167991
$ws.Range("D2").Value = 1.035493930877208
$ws.Range("E2").Value = 1.031277720443879
$ws.Range("F2").Value = 1.042354419174953
$ws.Range("I2").Value = 1.038203613957814
$ws.Range("J2").Value = 1.036890039641856
$ws.Range("K2").Value = 1.038290375760633
$ws.Range("L2").Value = 1.034086302204459
$ws.Range("M2").Value = 1.045131343117061
$ws.Range("N2").Value = 1.016232305288903

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032733174859362
$ws.Range("D3").Value = 1.036242997236936
$ws.Range("E3").Value = 1.032109421328874
$ws.Range("F3").Value = 1.044626218230895
$ws.Range("I3").Value = 1.038555124157658
$ws.Range("J3").Value = 1.037509017280812
$ws.Range("K3").Value = 1.038849008957261
$ws.Range("L3").Value = 1.034726471215681
$ws.Range("M3").Value = 1.047210133410427
$ws.Range("N3").Value = 1.016442318237395

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033365160267925
$ws.Range("D4").Value = 1.036727318410567
$ws.Range("E4").Value = 1.032647696103656
$ws.Range("F4").Value = 1.046090442384816
$ws.Range("I4").Value = 1.038780732827035
$ws.Range("J4").Value = 1.037908675624827
$ws.Range("K4").Value = 1.039209445185182
$ws.Range("L4").Value = 1.035140168100404
$ws.Range("M4").Value = 1.048549154384068
$ws.Range("N4").Value = 1.016577808094313

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033630783777393
$ws.Range("D5").Value = 1.036930838355711
$ws.Range("E5").Value = 1.032874013546064
$ws.Range("F5").Value = 1.046704652331262
$ws.Range("I5").Value = 1.038875139831378
$ws.Range("J5").Value = 1.038076486943559
$ws.Range("K5").Value = 1.039360725692556
$ws.Range("L5").Value = 1.03531395904579
$ws.Range("M5").Value = 1.049110651527169
$ws.Range("N5").Value = 1.016634671980089

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033675379471728
$ws.Range("D6").Value = 1.036965005080261
$ws.Range("E6").Value = 1.03291201482778
$ws.Range("F6").Value = 1.046807702712814
$ws.Range("I6").Value = 1.038890965506481
$ws.Range("J6").Value = 1.038104651212933
$ws.Range("K6").Value = 1.039386111916136
$ws.Range("L6").Value = 1.035343131875941
$ws.Range("M6").Value = 1.049204846592551
$ws.Range("N6").Value = 1.016644214053543

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033368709788092
$ws.Range("D7").Value = 1.03673003820076
$ws.Range("E7").Value = 1.03265072006212
$ws.Range("F7").Value = 1.046098654746279
$ws.Range("I7").Value = 1.038781996019592
$ws.Range("J7").Value = 1.037910918731701
$ws.Range("K7").Value = 1.039211467570789
$ws.Range("L7").Value = 1.035142490801938
$ws.Range("M7").Value = 1.048556662701635
$ws.Range("N7").Value = 1.01657856828922

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032086340734346
$ws.Range("D8").Value = 1.035747159626452
$ws.Range("E8").Value = 1.031558775581868
$ws.Range("F8").Value = 1.043123402216718
$ws.Range("I8").Value = 1.038322791124879
$ws.Range("J8").Value = 1.037099405168457
$ws.Range("K8").Value = 1.038479383952646
$ws.Range("L8").Value = 1.034302761712965
$ws.Range("M8").Value = 1.045835161479976
$ws.Range("N8").Value = 1.016303363700068

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02982450831325
$ws.Range("D9").Value = 1.034012281350796
$ws.Range("E9").Value = 1.029635428746176
$ws.Range("F9").Value = 1.03783479297128
$ws.Range("I9").Value = 1.037499411301794
$ws.Range("J9").Value = 1.035662759771511
$ws.Range("K9").Value = 1.037181356985607
$ws.Range("L9").Value = 1.032818909321772
$ws.Range("M9").Value = 1.040991397012442
$ws.Range("N9").Value = 1.015815319725696

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028315086536815
$ws.Range("D10").Value = 1.032853658021806
$ws.Range("E10").Value = 1.02835368321835
$ws.Range("F10").Value = 1.034276047029271
$ws.Range("I10").Value = 1.036940806550914
$ws.Range("J10").Value = 1.034700439535531
$ws.Range("K10").Value = 1.036310535899346
$ws.Range("L10").Value = 1.031826825492298
$ws.Range("M10").Value = 1.037727816413401
$ws.Range("N10").Value = 1.015487851311813

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027661103202078
$ws.Range("D11").Value = 1.032351460096395
$ws.Range("E11").Value = 1.0277987756494
$ws.Range("F11").Value = 1.032726749223003
$ws.Range("I11").Value = 1.036696596133865
$ws.Range("J11").Value = 1.034282643992903
$ws.Range("K11").Value = 1.035932142305793
$ws.Range("L11").Value = 1.031396551177457
$ws.Range("M11").Value = 1.036306029830958
$ws.Range("N11").Value = 1.015345548942699

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027418123373349
$ws.Range("D12").Value = 1.032164843679088
$ws.Range("E12").Value = 1.027592671571421
$ws.Range("F12").Value = 1.032149975653242
$ws.Range("I12").Value = 1.036605532530908
$ws.Range("J12").Value = 1.034127288390157
$ws.Range("K12").Value = 1.035791389484504
$ws.Range("L12").Value = 1.031236622219562
$ws.Range("M12").Value = 1.035776577779582
$ws.Range("N12").Value = 1.015292614912655

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027470246180289
$ws.Range("D13").Value = 1.032204877079679
$ws.Range("E13").Value = 1.027636881014091
$ws.Range("F13").Value = 1.032273754855575
$ws.Range("I13").Value = 1.036625082001427
$ws.Range("J13").Value = 1.034160620313975
$ws.Range("K13").Value = 1.035821590551198
$ws.Range("L13").Value = 1.031270932340437
$ws.Range("M13").Value = 1.035890208202245
$ws.Range("N13").Value = 1.015303972916273

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027641019676704
$ws.Range("D14").Value = 1.032336035920999
$ws.Range("E14").Value = 1.027781738757958
$ws.Range("F14").Value = 1.032679099615491
$ws.Range("I14").Value = 1.0366890760051
$ws.Range("J14").Value = 1.034269805682612
$ws.Range("K14").Value = 1.035920511737044
$ws.Range("L14").Value = 1.031383333569069
$ws.Range("M14").Value = 1.036262292671578
$ws.Range("N14").Value = 1.015341174964798

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027746230676335
$ws.Range("D15").Value = 1.032416836858779
$ws.Range("E15").Value = 1.027870992119802
$ws.Range("F15").Value = 1.03292867295622
$ws.Range("I15").Value = 1.03672845796493
$ws.Range("J15").Value = 1.034337056110758
$ws.Range("K15").Value = 1.03598143370542
$ws.Range("L15").Value = 1.031452573596814
$ws.Range("M15").Value = 1.036491367792816
$ws.Range("N15").Value = 1.015364086209563

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028358480779238
$ws.Range("D16").Value = 1.032886976434012
$ws.Range("E16").Value = 1.028390512539739
$ws.Range("F16").Value = 1.034378689528788
$ws.Range("I16").Value = 1.036956964646061
$ws.Range("J16").Value = 1.034728143826553
$ws.Range("K16").Value = 1.036335620612653
$ws.Range("L16").Value = 1.031855366607374
$ws.Range("M16").Value = 1.0378219903902
$ws.Range("N16").Value = 1.015497284738785

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028742422177564
$ws.Range("D17").Value = 1.033181745833857
$ws.Range("E17").Value = 1.02871641900989
$ws.Range("F17").Value = 1.035285984502261
$ws.Range("I17").Value = 1.037099674854803
$ws.Range("J17").Value = 1.034973165852463
$ws.Range("K17").Value = 1.036557437333023
$ws.Range("L17").Value = 1.032107840825111
$ws.Range("M17").Value = 1.03865431506755
$ws.Range("N17").Value = 1.015580700724194

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028966330857255
$ws.Range("D18").Value = 1.033353631096078
$ws.Range("E18").Value = 1.02890652415682
$ws.Range("F18").Value = 1.03581439280618
$ws.Range("I18").Value = 1.037182690601662
$ws.Range("J18").Value = 1.0351159766079
$ws.Range("K18").Value = 1.036686691822426
$ws.Range("L18").Value = 1.032255037761854
$ws.Range("M18").Value = 1.039138965585125
$ws.Range("N18").Value = 1.015629307034342

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029042671575709
$ws.Range("D19").Value = 1.03341223131051
$ws.Range("E19").Value = 1.028971346732494
$ws.Range("F19").Value = 1.035994431717515
$ws.Range("I19").Value = 1.037210958795898
$ws.Range("J19").Value = 1.035164653382948
$ws.Range("K19").Value = 1.036730742702575
$ws.Range("L19").Value = 1.032305216790912
$ws.Range("M19").Value = 1.039304079238099
$ws.Range("N19").Value = 1.015645872264043

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028701232828901
$ws.Range("D20").Value = 1.0331501249206
$ws.Range("E20").Value = 1.028681451364917
$ws.Range("F20").Value = 1.035188723548984
$ws.Range("I20").Value = 1.037084386662268
$ws.Range("J20").Value = 1.034946888331234
$ws.Range("K20").Value = 1.036533651694966
$ws.Range("L20").Value = 1.03208075967177
$ws.Range("M20").Value = 1.038565100684419
$ws.Range("N20").Value = 1.015571756032118

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027590732848049
$ws.Range("D21").Value = 1.032297415071155
$ws.Range("E21").Value = 1.027739081405795
$ws.Range("F21").Value = 1.032559771792923
$ws.Range("I21").Value = 1.036670241132733
$ws.Range("J21").Value = 1.034237657957873
$ws.Range("K21").Value = 1.035891387466432
$ws.Range("L21").Value = 1.031350237144367
$ws.Range("M21").Value = 1.036152760307186
$ws.Range("N21").Value = 1.015330222007333

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02689216167684
$ws.Range("D22").Value = 1.031760831856623
$ws.Range("E22").Value = 1.027146652321401
$ws.Range("F22").Value = 1.030899329051466
$ws.Range("I22").Value = 1.036407807738606
$ws.Range("J22").Value = 1.033790765283509
$ws.Range("K22").Value = 1.035486408895842
$ws.Range("L22").Value = 1.030890315080048
$ws.Range("M22").Value = 1.034628269178098
$ws.Range("N22").Value = 1.015177916562773

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027262521921483
$ws.Range("D23").Value = 1.032045328061884
$ws.Range("E23").Value = 1.027460703375682
$ws.Range("F23").Value = 1.031780288120976
$ws.Range("I23").Value = 1.036547123309097
$ws.Range("J23").Value = 1.034027764232184
$ws.Range("K23").Value = 1.035701206457459
$ws.Range("L23").Value = 1.031134187115228
$ws.Range("M23").Value = 1.035437179566533
$ws.Range("N23").Value = 1.015258698767416

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028719844639927
$ws.Range("D24").Value = 1.033164413201123
$ws.Range("E24").Value = 1.028697251709798
$ws.Range("F24").Value = 1.03523267406212
$ws.Range("I24").Value = 1.037091295432792
$ws.Range("J24").Value = 1.034958762341146
$ws.Range("K24").Value = 1.036544399794998
$ws.Range("L24").Value = 1.032092996686588
$ws.Range("M24").Value = 1.038605415388227
$ws.Range("N24").Value = 1.015575797904406

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030409509130143
$ws.Range("D25").Value = 1.034461142549751
$ws.Range("E25").Value = 1.030132570491965
$ws.Range("F25").Value = 1.039207680070532
$ws.Range("I25").Value = 1.037713971590971
$ws.Range("J25").Value = 1.036034964005709
$ws.Range("K25").Value = 1.037517885435008
$ws.Range("L25").Value = 1.033203018068721
$ws.Range("M25").Value = 1.04224953835432
$ws.Range("N25").Value = 1.015941860240291

Write-Output "Applied vm_pu.xlsx updates for case with 380 kV"